# Apply the daily cryptos-list refresh (prices + 1h % changes), matching the
# upstream GitHub Actions data pull. Cells are plain text (not numbers) in the
# source sheet, so numeric-looking Price values are forced to text via a
# temporary "@" number format (cleared again afterwards) to stop Excel from
# auto-converting typed entries like "678.29" into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textForceCells = @(
    'D5', 'D6', 'D8', 'D10', 'D13', 'D17', 'D19', 'D20', 'D22', 'D26', 'D33', 'D34', 'D36', 'D41', 'D43', 'D44', 'D45', 'D46', 'D47', 'D50', 'D51'
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '69.759.20'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '3.705.46'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '678.29'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '161.46'
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.497'
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").Value = '7.14'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  +1.83%  '
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '32.81'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '3.714.60'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '69.766.77'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("D17").Value = '16.08'
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '473.42'
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").Value = '9.81'
$ws.Range("E20").Value = '  -2.24%  '
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").Value = '80.53'
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("D23").Value = '3.852.73'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("D26").Value = '10.98'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").Value = '26.96'
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = '0.166'
$ws.Range("E34").Value = '  +3.62%  '
$ws.Range("D35").Value = '3.695.48'
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("D36").Value = '8.50'
$ws.Range("E36").Value = '  +3.64%  '
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '0.0906'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").Value = '166.90'
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").Value = '47.01'
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("D45").Value = '2.79'
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '28.50'
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D47").Value = '0.000281'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("D50").Value = '7.89'
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '0.267'
$ws.Range("E51").Value = '  +1.74%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
